# Update "想去人数" (column F) values on sheets "展览" and "全部类型"
# to match the latest scrape output (gh-pages update at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 331
    4  = 1219
    5  = 273
    7  = 60
    8  = 41
    10 = 3425
    11 = 123
    12 = 83
    16 = 584
    17 = 78
    18 = 708
    19 = 204
    20 = 115
    22 = 54
    24 = 2532
    25 = 5044
    29 = 1292
    30 = 279
    31 = 2219
    33 = 485
    35 = 97
    36 = 165
    39 = 784
    42 = 32
    43 = 467
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 331
    4  = 1219
    5  = 273
    7  = 60
    8  = 41
    10 = 3425
    11 = 123
    12 = 83
    17 = 584
    18 = 78
    19 = 708
    20 = 204
    21 = 115
    23 = 54
    25 = 2532
    26 = 5044
    30 = 1292
    31 = 279
    32 = 2219
    34 = 485
    36 = 97
    37 = 165
    40 = 784
    43 = 0
    44 = 467
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
